$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "2.90", "34.50", "1.00") keep their original text representation
# instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '48.005.23'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.491.93'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '318.17'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '104.84'
$ws.Range("E6").Value = '  -3.89%  '
$ws.Range("D7").Value = '0.519'
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  -2.10%  '
$ws.Range("D10").Value = '38.71'
$ws.Range("E10").Value = '  -2.98%  '
$ws.Range("D11").Value = '20.31'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '0.0799'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '7.06'
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").Value = '2.877.08'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").Value = '2.501.78'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '0.832'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '47.855.22'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").Value = '12.75'
$ws.Range("E19").Value = '  -3.32%  '
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").Value = '2.90'
$ws.Range("E20").Value = '  +7.63%  '
$ws.Range("D21").Value = '6.54'
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("D22").Value = '0.0₃0927'
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("D23").Value = '281.38'
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("D24").Value = '70.77'
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '25.65'
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  -7.38%  '
$ws.Range("D29").Value = '9.57'
$ws.Range("E29").Value = '  -5.15%  '
$ws.Range("D30").Value = '0.139'
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("D31").Value = '34.50'
$ws.Range("E31").Value = '  -2.59%  '
$ws.Range("D32").Value = '49.02'
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").Value = '19.21'
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").Value = '5.25'
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("D36").Value = '0.0769'
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").Value = '1.94'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").Value = '4.48'
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("D39").Value = '2.87'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("D42").Value = '119.13'
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").Value = '21.56'
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").Value = '0.0298'
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").Value = '1.986.59'
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("D46").Value = '3.10'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = '2.10'
$ws.Range("E47").Value = '  +5.79%  '
$ws.Range("D48").Value = '1.93'
$ws.Range("E48").Value = '  +4.28%  '
$ws.Range("D49").Value = '8.93'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").Value = '5.08'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").Value = '79.36'
$ws.Range("E51").Value = '  -0.70%  '

# Restore original (unset) formatting on column D so no style/number-format
# metadata is left behind.
$ws.Range("D2:D51").ClearFormats()
